$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 18.629453
$ws.Range("H2").Value = 55.888359
$ws.Range("I2").Value = 0.07116572597273459
$ws.Range("J2").Value = 0.07977938953593823
$ws.Range("M2").Value = 12.492041
$ws.Range("N2").Value = 37.476123
$ws.Range("O2").Value = 0.001845183434243001
$ws.Range("P2").Value = 0.001846965528678714
$ws.Range("Q2").Value = 232.719890683573
$ws.Range("R2").Value = 2094.479016152157
$ws.Range("S2").Value = 0.0001313138186507667
$ws.Range("T2").Value = 0.0001473497823719092
$ws.Range("G3").Value = 18.629453
$ws.Range("H3").Value = 55.888359
$ws.Range("I3").Value = 0.07116572597273459
$ws.Range("J3").Value = 0.07977938953593823
$ws.Range("O3").Value = 0.0009361563262120847
$ws.Range("P3").Value = 0.0009370604742490439
$ws.Range("Q3").Value = 118.070753213862
$ws.Range("R3").Value = 1062.636778924758
$ws.Range("S3").Value = 0.00006662224457885115
$ws.Range("T3").Value = 0.00007475811259384549
$ws.Range("G4").Value = 18.629453
$ws.Range("H4").Value = 55.888359
$ws.Range("I4").Value = 0.07116572597273459
$ws.Range("J4").Value = 0.07977938953593823
$ws.Range("M4").Value = 4688.500488333333
$ws.Range("N4").Value = 14065.501465
$ws.Range("O4").Value = 0.6925324238459419
$ws.Range("P4").Value = 0.6932012777691797
$ws.Range("Q4").Value = 87344.19948788288
$ws.Range("R4").Value = 786097.7953909459
$ws.Range("S4").Value = 0.04928457270265399
$ws.Range("T4").Value = 0.05530317476595751
$ws.Range("G5").Value = 18.629453
$ws.Range("H5").Value = 55.888359
$ws.Range("I5").Value = 0.07116572597273459
$ws.Range("J5").Value = 0.07977938953593823
$ws.Range("M5").Value = 19.5968845
$ws.Range("N5").Value = 39.193769
$ws.Range("O5").Value = 0.002894630800697294
$ws.Range("P5").Value = 0.001931617640437256
$ws.Range("Q5").Value = 365.0792387391785
$ws.Range("R5").Value = 2190.475432435071
$ws.Range("S5").Value = 0.0002059985023546609
$ws.Range("T5").Value = 0.0001541032761709337
$ws.Range("G6").Value = 18.629453
$ws.Range("H6").Value = 55.888359
$ws.Range("I6").Value = 0.07116572597273459
$ws.Range("J6").Value = 0.07977938953593823
$ws.Range("M6").Value = 2043.153564333333
$ws.Range("N6").Value = 6129.460693
$ws.Range("O6").Value = 0.3017916055929057
$ws.Range("P6").Value = 0.3020830785874553
$ws.Range("Q6").Value = 38062.83329853031
$ws.Range("R6").Value = 342565.4996867728
$ws.Range("S6").Value = 0.02147721870449632
$ws.Range("T6").Value = 0.02410000359884404
$ws.Range("I7").Value = 0.2779443552245922
$ws.Range("J7").Value = 0.3115858186182692
$ws.Range("M7").Value = 12.492041
$ws.Range("N7").Value = 37.476123
$ws.Range("O7").Value = 0.001845183434243001
$ws.Range("P7").Value = 0.001846965528678714
$ws.Range("Q7").Value = 908.9091564774462
$ws.Range("R7").Value = 8180.182408297017
$ws.Range("S7").Value = 0.0005128583199017694
$ws.Range("T7").Value = 0.0005754882662130814
$ws.Range("I8").Value = 0.2779443552245922
$ws.Range("J8").Value = 0.3115858186182692
$ws.Range("O8").Value = 0.0009361563262120847
$ws.Range("P8").Value = 0.0009370604742490439
$ws.Range("S8").Value = 0.0002601993664784409
$ws.Range("T8").Value = 0.0002919747549637119
$ws.Range("I9").Value = 0.2779443552245922
$ws.Range("J9").Value = 0.3115858186182692
$ws.Range("M9").Value = 4688.500488333333
$ws.Range("N9").Value = 14065.501465
$ws.Range("O9").Value = 0.6925324238459419
$ws.Range("P9").Value = 0.6932012777691797
$ws.Range("Q9").Value = 341130.8867778407
$ws.Range("R9").Value = 3070177.981000567
$ws.Range("S9").Value = 0.1924854780179843
$ws.Range("T9").Value = 0.2159916876009401
$ws.Range("I10").Value = 0.2779443552245922
$ws.Range("J10").Value = 0.3115858186182692
$ws.Range("M10").Value = 19.5968845
$ws.Range("N10").Value = 39.193769
$ws.Range("O10").Value = 0.002894630800697294
$ws.Range("P10").Value = 0.001931617640437256
$ws.Range("Q10").Value = 1425.850888616275
$ws.Range("R10").Value = 8555.10533169765
$ws.Range("S10").Value = 0.0008045462915130542
$ws.Range("T10").Value = 0.0006018646637531319
$ws.Range("I11").Value = 0.2779443552245922
$ws.Range("J11").Value = 0.3115858186182692
$ws.Range("M11").Value = 2043.153564333333
$ws.Range("N11").Value = 6129.460693
$ws.Range("O11").Value = 0.3017916055929057
$ws.Range("P11").Value = 0.3020830785874553
$ws.Range("Q11").Value = 148657.9320954916
$ws.Range("R11").Value = 1337921.388859425
$ws.Range("S11").Value = 0.0838812732287146
$ws.Range("T11").Value = 0.0941248033323992
$ws.Range("G12").Value = 36.272704
$ws.Range("H12").Value = 108.818112
$ws.Range("I12").Value = 0.1385640959589159
$ws.Range("J12").Value = 0.1553354348051864
$ws.Range("M12").Value = 12.492041
$ws.Range("N12").Value = 37.476123
$ws.Range("O12").Value = 0.001845183434243001
$ws.Range("P12").Value = 0.001846965528678714
$ws.Range("Q12").Value = 453.120105548864
$ws.Range("R12").Value = 4078.080949939776
$ws.Range("S12").Value = 0.0002556761744442491
$ws.Range("T12").Value = 0.000286899193467499
$ws.Range("G13").Value = 36.272704
$ws.Range("H13").Value = 108.818112
$ws.Range("I13").Value = 0.1385640959589159
$ws.Range("J13").Value = 0.1553354348051864
$ws.Range("O13").Value = 0.0009361563262120847
$ws.Range("P13").Value = 0.0009370604742490439
$ws.Range("Q13").Value = 229.891102137216
$ws.Range("R13").Value = 2069.019919234944
$ws.Range("S13").Value = 0.0001297176550177975
$ws.Range("T13").Value = 0.0001455586962062294
$ws.Range("G14").Value = 36.272704
$ws.Range("H14").Value = 108.818112
$ws.Range("I14").Value = 0.1385640959589159
$ws.Range("J14").Value = 0.1553354348051864
$ws.Range("M14").Value = 4688.500488333333
$ws.Range("N14").Value = 14065.501465
$ws.Range("O14").Value = 0.6925324238459419
$ws.Range("P14").Value = 0.6932012777691797
$ws.Range("Q14").Value = 170064.5904171704
$ws.Range("R14").Value = 1530581.313754534
$ws.Range("S14").Value = 0.0959601292324497
$ws.Range("T14").Value = 0.1076787218897863
$ws.Range("G15").Value = 36.272704
$ws.Range("H15").Value = 108.818112
$ws.Range("I15").Value = 0.1385640959589159
$ws.Range("J15").Value = 0.1553354348051864
$ws.Range("M15").Value = 19.5968845
$ws.Range("N15").Value = 39.193769
$ws.Range("O15").Value = 0.002894630800697294
$ws.Range("P15").Value = 0.001931617640437256
$ws.Range("Q15").Value = 710.8319907906879
$ws.Range("R15").Value = 4264.991944744128
$ws.Range("S15").Value = 0.0004010919000334533
$ws.Range("T15").Value = 0.0003000486660546894
$ws.Range("G16").Value = 36.272704
$ws.Range("H16").Value = 108.818112
$ws.Range("I16").Value = 0.1385640959589159
$ws.Range("J16").Value = 0.1553354348051864
$ws.Range("M16").Value = 2043.153564333333
$ws.Range("N16").Value = 6129.460693
$ws.Range("O16").Value = 0.3017916055929057
$ws.Range("P16").Value = 0.3020830785874553
$ws.Range("Q16").Value = 74110.70446560795
$ws.Range("R16").Value = 666996.3401904716
$ws.Range("S16").Value = 0.04181748099697068
$ws.Range("T16").Value = 0.04692420635967167
$ws.Range("G17").Value = 84.7905925
$ws.Range("H17").Value = 169.581185
$ws.Range("I17").Value = 0.3239055956672912
$ws.Range("J17").Value = 0.2420733701642771
$ws.Range("M17").Value = 12.492041
$ws.Range("N17").Value = 37.476123
$ws.Range("O17").Value = 0.001845183434243001
$ws.Range("P17").Value = 0.001846965528678714
$ws.Range("Q17").Value = 1059.207557924293
$ws.Range("R17").Value = 6355.245347545755
$ws.Range("S17").Value = 0.0005976652393838971
$ws.Range("T17").Value = 0.0004471011701045019
$ws.Range("G18").Value = 84.7905925
$ws.Range("H18").Value = 169.581185
$ws.Range("I18").Value = 0.3239055956672912
$ws.Range("J18").Value = 0.2420733701642771
$ws.Range("O18").Value = 0.0009361563262120847
$ws.Range("P18").Value = 0.0009370604742490439
$ws.Range("Q18").Value = 537.390395838495
$ws.Range("R18").Value = 3224.34237503097
$ws.Range("S18").Value = 0.0003032262724794283
$ws.Range("T18").Value = 0.0002268373870492018
$ws.Range("G19").Value = 84.7905925
$ws.Range("H19").Value = 169.581185
$ws.Range("I19").Value = 0.3239055956672912
$ws.Range("J19").Value = 0.2420733701642771
$ws.Range("M19").Value = 4688.500488333333
$ws.Range("N19").Value = 14065.501465
$ws.Range("O19").Value = 0.6925324238459419
$ws.Range("P19").Value = 0.6932012777691797
$ws.Range("Q19").Value = 397540.7343423227
$ws.Range("R19").Value = 2385244.406053936
$ws.Range("S19").Value = 0.2243151272647328
$ws.Range("T19").Value = 0.1678055695117685
$ws.Range("G20").Value = 84.7905925
$ws.Range("H20").Value = 169.581185
$ws.Range("I20").Value = 0.3239055956672912
$ws.Range("J20").Value = 0.2420733701642771
$ws.Range("M20").Value = 19.5968845
$ws.Range("N20").Value = 39.193769
$ws.Range("O20").Value = 0.002894630800697294
$ws.Range("P20").Value = 0.001931617640437256
$ws.Range("Q20").Value = 1661.631447909066
$ws.Range("R20").Value = 6646.525791636264
$ws.Range("S20").Value = 0.0009375871137367449
$ws.Range("T20").Value = 0.0004675931920894152
$ws.Range("G21").Value = 84.7905925
$ws.Range("H21").Value = 169.581185
$ws.Range("I21").Value = 0.3239055956672912
$ws.Range("J21").Value = 0.2420733701642771
$ws.Range("M21").Value = 2043.153564333333
$ws.Range("N21").Value = 6129.460693
$ws.Range("O21").Value = 0.3017916055929057
$ws.Range("P21").Value = 0.3020830785874553
$ws.Range("Q21").Value = 173240.2012883102
$ws.Range("R21").Value = 1039441.207729861
$ws.Range("S21").Value = 0.09775198977695831
$ws.Range("T21").Value = 0.07312626890326546
$ws.Range("G22").Value = 49.32382433333333
$ws.Range("H22").Value = 147.971473
$ws.Range("I22").Value = 0.1884202271764661
$ws.Range("J22").Value = 0.2112259868763291
$ws.Range("M22").Value = 12.492041
$ws.Range("N22").Value = 37.476123
$ws.Range("O22").Value = 0.001845183434243001
$ws.Range("P22").Value = 0.001846965528678714
$ws.Range("Q22").Value = 616.1552358487977
$ws.Range("R22").Value = 5545.39712263918
$ws.Range("S22").Value = 0.000347669881862318
$ws.Range("T22").Value = 0.0003901271165217221
$ws.Range("G23").Value = 49.32382433333333
$ws.Range("H23").Value = 147.971473
$ws.Range("I23").Value = 0.1884202271764661
$ws.Range("J23").Value = 0.2112259868763291
$ws.Range("O23").Value = 0.0009361563262120847
$ws.Range("P23").Value = 0.0009370604742490439
$ws.Range("Q23").Value = 312.607197346314
$ws.Range("R23").Value = 2813.464776116826
$ws.Range("S23").Value = 0.0001763907876575669
$ws.Range("T23").Value = 0.0001979315234360552
$ws.Range("G24").Value = 49.32382433333333
$ws.Range("H24").Value = 147.971473
$ws.Range("I24").Value = 0.1884202271764661
$ws.Range("J24").Value = 0.2112259868763291
$ws.Range("M24").Value = 4688.500488333333
$ws.Range("N24").Value = 14065.501465
$ws.Range("O24").Value = 0.6925324238459419
$ws.Range("P24").Value = 0.6932012777691797
$ws.Range("Q24").Value = 231254.7744733009
$ws.Range("R24").Value = 2081292.970259708
$ws.Range("S24").Value = 0.1304871166281211
$ws.Range("T24").Value = 0.1464221240007273
$ws.Range("G25").Value = 49.32382433333333
$ws.Range("H25").Value = 147.971473
$ws.Range("I25").Value = 0.1884202271764661
$ws.Range("J25").Value = 0.2112259868763291
$ws.Range("M25").Value = 19.5968845
$ws.Range("N25").Value = 39.193769
$ws.Range("O25").Value = 0.002894630800697294
$ws.Range("P25").Value = 0.001931617640437256
$ws.Range("Q25").Value = 966.5932885586227
$ws.Range("R25").Value = 5799.559731351736
$ws.Range("S25").Value = 0.00054540699305938
$ws.Range("T25").Value = 0.0004080078423690855
$ws.Range("G26").Value = 49.32382433333333
$ws.Range("H26").Value = 147.971473
$ws.Range("I26").Value = 0.1884202271764661
$ws.Range("J26").Value = 0.2112259868763291
$ws.Range("M26").Value = 2043.153564333333
$ws.Range("N26").Value = 6129.460693
$ws.Range("O26").Value = 0.3017916055929057
$ws.Range("P26").Value = 0.3020830785874553
$ws.Range("Q26").Value = 100776.1474932012
$ws.Range("R26").Value = 906985.3274388108
$ws.Range("S26").Value = 0.05686364288576574
$ws.Range("T26").Value = 0.06380779639327491
